$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Register")

# Update the two registered email addresses (shown via hyperlink-styled cells)
$ws.Range("C2").Value = "bt6g7h3b4@gmail.com"
$ws.Range("C3").Value = "t5h8d4r6v8@gmail.com"

# Move the active selection to D5 as recorded in the saved view state
$ws.Activate()
$ws.Range("D5").Select()
